# Commit: "Changed one slide to remove ugly arrow"
#
# The "Straight Connector 6" shape on the "1. Design: then vs Now" slide
# (slide 4 - the one with "Iteration 1: The loading of different UI's was
# handled in different UIElements") is a stray thin line left over from
# editing that renders like an ugly little arrow/stroke poking out of the
# corner of the slide. Remove it, leaving the title and content placeholder
# untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Straight Connector 6") {
        $shp.Delete()
    }
}
